$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 10840
$ws1.Range("F15").Value = 9058
$ws1.Range("F17").Value = 736
$ws1.Range("F20").Value = 3370

# Sheet "全部类型" (All Types) - update "想去人数" (want-to-go count) column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 10841
$ws4.Range("F18").Value = 9058
$ws4.Range("F20").Value = 736
$ws4.Range("F23").Value = 3370
